$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$oldGuid = "12ff1f41-8210-45e4-bddd-5ef673472969"
$newGuid = "a430babd-8ead-4cf0-8f75-f101c70c5bd4"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af61115601566d82751cf29e7b82e1dacac981a2/e2e/" + $oldGuid + ".md"

$hlColor = 15570276   # BGR-packed int for FF6495ED

# ---------------------------------------------------------------------------
# Sheet1 "Overview"
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = $newGuid + ".md"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $baseUrl, "", "", "e2e\" + $newGuid + ".md") | Out-Null
$ws1.Range("B2").Font.Color = $hlColor
$ws1.Range("B2").Font.Underline = $true

$ws1.Range("G2").Value = "2016-08-26 09:05:28"

# ---------------------------------------------------------------------------
# Sheet2 "zh-cn"
# ---------------------------------------------------------------------------
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $baseUrl, "", "", $newGuid + ".md") | Out-Null
$ws2.Range("A2").Font.Color = $hlColor
$ws2.Range("A2").Font.Underline = $true

$ws2.Range("G2").Value = $newGuid + ".50ab5d8de95dc41ea51f31b2f89d08a60ebd2e27.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-26 09:05:24"
$ws2.Range("I2").Value = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws2.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# Sheet3 "de-de"
# ---------------------------------------------------------------------------
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $baseUrl, "", "", $newGuid + ".md") | Out-Null
$ws3.Range("A2").Font.Color = $hlColor
$ws3.Range("A2").Font.Underline = $true

$ws3.Range("G2").Value = $newGuid + ".50ab5d8de95dc41ea51f31b2f89d08a60ebd2e27.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-26 09:05:28"
$ws3.Range("I2").Value = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws3.Columns.Item(10).ColumnWidth = 21.7054770333426
